$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (co2)
$ws.Range("C2").Value = 532.2106212425854
$ws.Range("D2").Value = 126.9954844019068
$ws.Range("G2").Value = 485
$ws.Range("H2").Value = 584

# Row 3 (humidity)
$ws.Range("C3").Value = 36.51327065315416
$ws.Range("D3").Value = 6.585029651722805
$ws.Range("F3").Value = 31.67
$ws.Range("G3").Value = 36.67
$ws.Range("H3").Value = 40.85

# Row 4 (pm25)
$ws.Range("C4").Value = 2.051220549113259
$ws.Range("D4").Value = 2.674862232545927
$ws.Range("F4").Value = 0.67
$ws.Range("G4").Value = 1.31
$ws.Range("H4").Value = 2.47

# Row 5 (pressure)
$ws.Range("C5").Value = 323.0362109370496
$ws.Range("D5").Value = 9.386588224888023
$ws.Range("F5").Value = 317.87
$ws.Range("G5").Value = 323.81
$ws.Range("H5").Value = 329.8

# Row 6 (temperature)
$ws.Range("C6").Value = 22.58150089553854
$ws.Range("D6").Value = 3.266941800073979
$ws.Range("F6").Value = 20.42
$ws.Range("G6").Value = 22.28
$ws.Range("H6").Value = 24.51
$ws.Range("I6").Value = 45.32

# Row 7 (rssi)
$ws.Range("C7").Value = -76.69644376203388
$ws.Range("D7").Value = 22.87019877802588

# Row 8 (snr)
$ws.Range("C8").Value = 7.415354158273629
$ws.Range("D8").Value = 7.104336290239358

# Row 9 (SF)
$ws.Range("C9").Value = 9.321388809101659
$ws.Range("D9").Value = 1.685678676124756

# Row 10 (frequency)
$ws.Range("C10").Value = 867.8302548663114
$ws.Range("D10").Value = 0.4616805344479739

# Row 11 (toa)
$ws.Range("C11").Value = 0.555527425568039
$ws.Range("D11").Value = 0.5888062665095325

# Row 12 (distance)
$ws.Range("C12").Value = 22.76446408797683
$ws.Range("D12").Value = 12.29682031832352

# Row 13 (c_walls)
$ws.Range("C13").Value = 0.6736605395087145
$ws.Range("D13").Value = 0.7508222926507676

# Row 14 (w_walls)
$ws.Range("C14").Value = 1.830988090064159
$ws.Range("D14").Value = 1.666477230061925

# Row 15 (exp_pl)
$ws.Range("C15").Value = 93.95644376203364
$ws.Range("D15").Value = 22.87019877802589

# Row 16 (n_power)
$ws.Range("C16").Value = -85.81856420269712
$ws.Range("D16").Value = 20.45211142633633
$ws.Range("F16").Value = -102.3377954106368
$ws.Range("H16").Value = -70.07382219273629

# Row 17 (esp)
$ws.Range("C17").Value = -78.40321004442367
$ws.Range("D17").Value = 25.42481706855655
$ws.Range("F17").Value = -93.49305820175223
$ws.Range("G17").Value = -74.16954289279533
$ws.Range("H17").Value = -60.25410721860875
